$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (scheduled GitHub Actions scrape).
# The Price/Volume(1h) columns store plain text (e.g. "297.14", "-0.38%"),
# not numbers, so each cell is switched to text format before the new
# literal is written (otherwise Excel would silently reinterpret
# "297.85" as a number or "-0.11%" as a percentage). The style is then
# reset to "Normal" so only the underlying value changes, matching the
# original (unstyled) data cells.
$updates = @{
    "D2" = "297.85"
    "E2" = "-0.11%"
    "D3" = "31.38"
    "E3" = "-0.26%"
    "D4" = "5.095"
    "E4" = "-0.96%"
    "D5" = "0.08000"
    "E5" = "9.29%"
    "D6" = "2.457"
    "E6" = "33.20%"
    "D7" = "7.795"
    "E7" = "0.42%"
    "D8" = "3.799"
    "E8" = "1.54%"
    "D9" = "0.9224"
    "E9" = "-0.36%"
    "D10" = "0.1728"
    "E10" = "3.24%"
    "D11" = "0.07314"
    "E11" = "3.52%"
    "D12" = "0.09037"
    "E12" = "13.36%"
    "D13" = "0.03028"
    "E13" = "0.87%"
    "D14" = "0.09984"
    "E14" = "0.58%"
    "D15" = "0.001490"
    "E15" = "-1.14%"
    "D16" = "0.005955"
    "E16" = "-3.18%"
    "D17" = "3.505"
    "E17" = "1.40%"
    "D18" = "2.248"
    "E18" = "1.28%"
    "E19" = "1.83%"
    "D20" = "0.1337"
    "E20" = "1.55%"
    "D21" = "4.585"
    "E21" = "0.49%"
    "D23" = "0.04589"
    "E23" = "-1.13%"
    "D24" = "0.001244"
    "E24" = "2.41%"
    "D25" = "0.004437"
    "E25" = "-6.24%"
    "D26" = "0.0001200"
    "E26" = "-7.67%"
    "D27" = "0.0003427"
    "E27" = "82.86%"
    "D39" = "0.01786"
    "E39" = "4.37%"
    "D40" = "0.04463"
    "E40" = "-0.09%"
    "D41" = "0.006984"
    "E41" = "-2.33%"
    "D42" = "0.1343"
    "E42" = "0.85%"
    "D43" = "0.002129"
    "E43" = "-0.93%"
    "E44" = "-6.55%"
    "D45" = "0.00006566"
    "E45" = "5.46%"
    "E46" = "-0.01%"
    "E48" = "11.00%"
    "D49" = "0.00002100"
    "E49" = "-0.01%"
    "D50" = "0.0002000"
    "E50" = "0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
